# Update "F" column (想去人数 / want-to-go count) values across sheets
# to reflect the regenerated data snapshot (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Cells.Item(6, 6).Value = 953
$ws1.Cells.Item(9, 6).Value = 2241
$ws1.Cells.Item(10, 6).Value = 640
$ws1.Cells.Item(11, 6).Value = 307
$ws1.Cells.Item(13, 6).Value = 1138
$ws1.Cells.Item(14, 6).Value = 191
$ws1.Cells.Item(15, 6).Value = 2260
$ws1.Cells.Item(16, 6).Value = 710
$ws1.Cells.Item(17, 6).Value = 14481
$ws1.Cells.Item(19, 6).Value = 1379
$ws1.Cells.Item(20, 6).Value = 341
$ws1.Cells.Item(21, 6).Value = 572
$ws1.Cells.Item(22, 6).Value = 147
$ws1.Cells.Item(23, 6).Value = 310
$ws1.Cells.Item(24, 6).Value = 153
$ws1.Cells.Item(25, 6).Value = 108
$ws1.Cells.Item(26, 6).Value = 47
$ws1.Cells.Item(29, 6).Value = 20
$ws1.Cells.Item(30, 6).Value = 9

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(6, 6).Value = 15
$ws2.Cells.Item(16, 6).Value = 5

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Cells.Item(2, 6).Value = 5745
$ws3.Cells.Item(4, 6).Value = 482

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Cells.Item(4, 6).Value = 482
$ws4.Cells.Item(8, 6).Value = 953
$ws4.Cells.Item(12, 6).Value = 2241
$ws4.Cells.Item(13, 6).Value = 640
$ws4.Cells.Item(14, 6).Value = 307
$ws4.Cells.Item(15, 6).Value = 15
$ws4.Cells.Item(18, 6).Value = 1138
$ws4.Cells.Item(20, 6).Value = 191
$ws4.Cells.Item(23, 6).Value = 2260
$ws4.Cells.Item(24, 6).Value = 710
$ws4.Cells.Item(28, 6).Value = 1379
$ws4.Cells.Item(29, 6).Value = 341
$ws4.Cells.Item(30, 6).Value = 572
$ws4.Cells.Item(31, 6).Value = 147
$ws4.Cells.Item(32, 6).Value = 310
$ws4.Cells.Item(33, 6).Value = 153
$ws4.Cells.Item(34, 6).Value = 108
$ws4.Cells.Item(36, 6).Value = 47
$ws4.Cells.Item(40, 6).Value = 5
$ws4.Cells.Item(42, 6).Value = 0
$ws4.Cells.Item(44, 6).Value = 9
